# Adapt column header formatting to respective input file names:
#   <field>_old -> <field>_FV2304
#   <field>_new -> <field>_FV2310
# Then turn the header row + data range into an Excel Table ("Table1"),
# and freeze the header row (row 1) so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header names, column A (1) .. U (21), in order.
$headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304",
    "diff",
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn A1:U59 into a real Excel Table ("Table1") with the renamed headers,
# an autofilter and the default banded-rows table style.
$rng = $ws.Range("A1:U59")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)

# Freeze the header row (row 1) by splitting/freezing above row 2.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Restore the selection to A1 (top of sheet) once the freeze is in place.
$ws.Range("A1").Select()
